$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries legacy sheet protection (password-hash "D382"). We don't
# know the clear-text password, so unprotect without one (Excel allows this
# programmatically) to unlock the cells for editing, matching what the
# original author's tool must have done to update these figures.
$ws.Unprotect()

# Update the "as of" date in the confidential notice (A9).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) / Percent Change (E) figures for each sector row.
$ws.Range("D2").Value = 0.2563861054688267
$ws.Range("E2").Value = 0.008758958025253172

$ws.Range("D3").Value = 0.2549651135532997
$ws.Range("E3").Value = 0.001847940865892461

$ws.Range("D4").Value = 0.2448752386022326
$ws.Range("E4").Value = -0.003055475985868417

$ws.Range("D5").Value = 0.2437735423756411
$ws.Range("E5").Value = -0.007912829160721335

$ws.Range("E6").Value = 0.00003968678296639716
